$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Games")
$ws.Range("C138").Copy()
$ws.Range("C140").PasteSpecial(-4122)
